$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price (column D) updates ---
$ws.Range("D2").Value = "29.358.19"
$ws.Range("D3").Value = "1.863.93"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.10"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07902"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3120"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.37"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07794"
$ws.Range("D12").Value = "1.879.32"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.146"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "92.35"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6978"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.534"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008556"
$ws.Range("D18").Value = "29.367.03"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "248.64"
$ws.Range("D20").Value = "2.121.11"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.000"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.583"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1537"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.971"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "160.56"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.73"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.586"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.294"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.244"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.199"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05247"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.885"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7556"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.180"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.711"
$ws.Range("D38").Value = "1.277.03"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01866"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.748"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8969"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "109.76"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.951"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "70.12"
$ws.Range("D46").Value = "2.021.25"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.579"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.793"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.5172"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4288"

# --- Volume(1h) (column E) updates ---
$ws.Range("E2").Value = "  -0.01%  "
$ws.Range("E3").Value = "  -0.81%  "
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("E5").Value = "  +0.03%  "
$ws.Range("E6").Value = "  -2.40%  "
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("E8").Value = "  -0.78%  "
$ws.Range("E9").Value = "  -0.71%  "
$ws.Range("E10").Value = "  -1.99%  "
$ws.Range("E11").Value = "  -4.34%  "
$ws.Range("E12").Value = "  -0.06%  "
$ws.Range("E13").Value = "  -1.60%  "
$ws.Range("E14").Value = "  -2.53%  "
$ws.Range("E15").Value = "  -1.69%  "
$ws.Range("E16").Value = "  +2.10%  "
$ws.Range("E17").Value = "  +1.29%  "
$ws.Range("E18").Value = "  +0.01%  "
$ws.Range("E19").Value = "  -0.20%  "
$ws.Range("E20").Value = "  -0.71%  "
$ws.Range("E21").Value = "  -2.36%  "
$ws.Range("E22").Value = "  -0.14%  "
$ws.Range("E23").Value = "  -2.05%  "
$ws.Range("E24").Value = "  -0.35%  "
$ws.Range("E25").Value = "  -3.41%  "
$ws.Range("E26").Value = "  -1.09%  "
$ws.Range("E27").Value = "  -1.25%  "
$ws.Range("E28").Value = "  -0.81%  "
$ws.Range("E29").Value = "  +5.44%  "
$ws.Range("E31").Value = "  -0.92%  "
$ws.Range("E32").Value = "  -1.78%  "
$ws.Range("E33").Value = "  -1.60%  "
$ws.Range("E34").Value = "  -2.74%  "
$ws.Range("E35").Value = "  -0.04%  "
$ws.Range("E36").Value = "  +0.08%  "
$ws.Range("E37").Value = "  +0.36%  "
$ws.Range("E38").Value = "  +0.73%  "
$ws.Range("E39").Value = "  -1.06%  "
$ws.Range("E40").Value = "  -0.57%  "
$ws.Range("E41").Value = "  -0.99%  "
$ws.Range("E42").Value = "  -3.11%  "
$ws.Range("E43").Value = "  -7.66%  "
$ws.Range("E44").Value = "  -5.90%  "
$ws.Range("E45").Value = "  -0.08%  "
$ws.Range("E46").Value = "  +0.00%  "
$ws.Range("E47").Value = "  -3.48%  "
$ws.Range("E48").Value = "  +0.81%  "
$ws.Range("E49").Value = "  -0.39%  "
$ws.Range("E50").Value = "  -0.52%  "
$ws.Range("E51").Value = "  -1.17%  "
